$wb = $excel.ActiveWorkbook

# Overview sheet: status text changes from "Ready for handoff" to
# "Handed back: in sync with en-US" for the two tracked files.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: mark the two tracked files as handed back, filling in the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("E2").Value = "2ecd04fc-f6f0-409d-8003-2c00f8de9488.md"
$wsZh.Range("F2").Value = "2ecd04fc-f6f0-409d-8003-2c00f8de9488.4684df10d29f7b7b24e648c49423aa14716b8531.zh-cn.xlf"
$wsZh.Range("E2").Style = "HyperLink"
$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Range("G2").Value = "2016-03-09 18:40:52"

$wsZh.Range("B3").Value = "Handed back: in sync with en-US"
$wsZh.Range("E3").Value = "6c476987-2e90-41aa-b3d3-52976c12f200.md"
$wsZh.Range("F3").Value = "6c476987-2e90-41aa-b3d3-52976c12f200.b356f391b1146cbf5449f0a559aea847d5c292a5.zh-cn.xlf"
$wsZh.Range("E3").Style = "HyperLink"
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Range("G3").Value = "2016-03-09 18:40:52"

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/de2b7a9e6b5365b88c29d26ac37562dbbb0db92f/e2e/2ecd04fc-f6f0-409d-8003-2c00f8de9488.md", "", "", "2ecd04fc-f6f0-409d-8003-2c00f8de9488.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0a734439805983efa74024581806ede78b4e674/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2ecd04fc-f6f0-409d-8003-2c00f8de9488.4684df10d29f7b7b24e648c49423aa14716b8531.zh-cn.xlf", "", "", "2ecd04fc-f6f0-409d-8003-2c00f8de9488.4684df10d29f7b7b24e648c49423aa14716b8531.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/de2b7a9e6b5365b88c29d26ac37562dbbb0db92f/e2e/6c476987-2e90-41aa-b3d3-52976c12f200.md", "", "", "6c476987-2e90-41aa-b3d3-52976c12f200.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0a734439805983efa74024581806ede78b4e674/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6c476987-2e90-41aa-b3d3-52976c12f200.b356f391b1146cbf5449f0a559aea847d5c292a5.zh-cn.xlf", "", "", "6c476987-2e90-41aa-b3d3-52976c12f200.b356f391b1146cbf5449f0a559aea847d5c292a5.zh-cn.xlf")

# de-de sheet: same handback updates.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("E2").Value = "2ecd04fc-f6f0-409d-8003-2c00f8de9488.md"
$wsDe.Range("F2").Value = "2ecd04fc-f6f0-409d-8003-2c00f8de9488.4684df10d29f7b7b24e648c49423aa14716b8531.de-de.xlf"
$wsDe.Range("E2").Style = "HyperLink"
$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Range("G2").Value = "2016-03-09 18:41:01"

$wsDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDe.Range("E3").Value = "6c476987-2e90-41aa-b3d3-52976c12f200.md"
$wsDe.Range("F3").Value = "6c476987-2e90-41aa-b3d3-52976c12f200.b356f391b1146cbf5449f0a559aea847d5c292a5.de-de.xlf"
$wsDe.Range("E3").Style = "HyperLink"
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Range("G3").Value = "2016-03-09 18:41:01"

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/de2b7a9e6b5365b88c29d26ac37562dbbb0db92f/e2e/2ecd04fc-f6f0-409d-8003-2c00f8de9488.md", "", "", "2ecd04fc-f6f0-409d-8003-2c00f8de9488.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b02577c86e3ea7e54c59b28e30b63c4859b5e83/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2ecd04fc-f6f0-409d-8003-2c00f8de9488.4684df10d29f7b7b24e648c49423aa14716b8531.de-de.xlf", "", "", "2ecd04fc-f6f0-409d-8003-2c00f8de9488.4684df10d29f7b7b24e648c49423aa14716b8531.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/de2b7a9e6b5365b88c29d26ac37562dbbb0db92f/e2e/6c476987-2e90-41aa-b3d3-52976c12f200.md", "", "", "6c476987-2e90-41aa-b3d3-52976c12f200.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b02577c86e3ea7e54c59b28e30b63c4859b5e83/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6c476987-2e90-41aa-b3d3-52976c12f200.b356f391b1146cbf5449f0a559aea847d5c292a5.de-de.xlf", "", "", "6c476987-2e90-41aa-b3d3-52976c12f200.b356f391b1146cbf5449f0a559aea847d5c292a5.de-de.xlf")
